$d = $word.ActiveDocument

# 1) The paragraph starting with "c)" had a stale <w:lastRenderedPageBreak/>
#    layout hint left over from a previous render pass. The run holding
#    "c)" is edited (identity replace) so the stale page-break marker gets
#    dropped, the way Word drops it once that text is touched again.
$d.Content.Find.Execute("c)", $true, $false, $false, $false, $false, $true, 1, $false, "c)", 2)

# 2) The stale "_GoBack" bookmark pair (marking the previous session's last
#    edit location, sitting between the "animaton" run and the following
#    space run) is removed outright via the Bookmarks collection, which
#    deletes just the bookmarkStart/bookmarkEnd pair without disturbing the
#    surrounding runs.
$d.Bookmarks("_GoBack").Delete()
